$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The dataset originally enumerated all 3x3 sending/target cluster combinations
# (including self-pairs). The refreshed TPM run drops the 3 self-pair rows
# (old rows 2, 6 and 10) and recomputes the NATMI metrics for the remaining
# 6 sending->target combinations, which shift up into rows 2-7.
# Delete the old trailing self-pair row (MuSCs -> MuSCs) and the two rows
# that already got superseded by recomputed values further up the sheet.
$ws.Range("A8:T10").EntireRow.Delete()

# Row 2: updated NATMI metrics for this sending/target cluster pair
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.535365
$ws.Range("H2").Value = 1.606095
$ws.Range("I2").Value = 0.1618182173563651
$ws.Range("J2").Value = 0.1618182173563651
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.219191666666667
$ws.Range("N2").Value = 3.657575
$ws.Range("O2").Value = 0.8635205180780776
$ws.Range("P2").Value = 0.8635205180780777
$ws.Range("Q2").Value = 0.652712546625
$ws.Range("R2").Value = 5.874412919625001
$ws.Range("S2").Value = 0.1397333508860394
$ws.Range("T2").Value = 0.1397333508860394

# Row 3: updated NATMI metrics for this sending/target cluster pair
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.535365
$ws.Range("H3").Value = 1.606095
$ws.Range("I3").Value = 0.1618182173563651
$ws.Range("J3").Value = 0.1618182173563651
$ws.Range("M3").Value = 0.1926933333333333
$ws.Range("N3").Value = 0.5780799999999999
$ws.Range("O3").Value = 0.1364794819219223
$ws.Range("P3").Value = 0.1364794819219224
$ws.Range("Q3").Value = 0.1031612664
$ws.Range("R3").Value = 0.9284513976
$ws.Range("S3").Value = 0.02208486647032573
$ws.Range("T3").Value = 0.02208486647032573

# Row 4: updated NATMI metrics for this sending/target cluster pair
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "FAPs"
$ws.Range("G4").Value = 2.059271
$ws.Range("H4").Value = 6.177813
$ws.Range("I4").Value = 0.6224306076670297
$ws.Range("J4").Value = 0.6224306076670296
$ws.Range("M4").Value = 1.219191666666667
$ws.Range("N4").Value = 3.657575
$ws.Range("O4").Value = 0.8635205180780776
$ws.Range("P4").Value = 0.8635205180780777
$ws.Range("Q4").Value = 2.510646042608334
$ws.Range("R4").Value = 22.595814383475
$ws.Range("S4").Value = 0.5374816008002862
$ws.Range("T4").Value = 0.5374816008002862

# Row 5: updated NATMI metrics for this sending/target cluster pair
$ws.Range("D5").Value = "MuSCs"
$ws.Range("I5").Value = 0.6224306076670297
$ws.Range("J5").Value = 0.6224306076670296
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.1926933333333333
$ws.Range("N5").Value = 0.5780799999999999
$ws.Range("O5").Value = 0.1364794819219223
$ws.Range("P5").Value = 0.1364794819219224
$ws.Range("Q5").Value = 0.3968077932266666
$ws.Range("R5").Value = 3.57127013904
$ws.Range("S5").Value = 0.08494900686674352
$ws.Range("T5").Value = 0.08494900686674352

# Row 6: updated NATMI metrics for this sending/target cluster pair
$ws.Range("A6").Value = "MuSCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.7137986666666666
$ws.Range("H6").Value = 2.141396
$ws.Range("I6").Value = 0.2157511749766052
$ws.Range("J6").Value = 0.2157511749766052
$ws.Range("O6").Value = 0.8635205180780776
$ws.Range("P6").Value = 0.8635205180780777
$ws.Range("Q6").Value = 0.8702573860777777
$ws.Range("R6").Value = 7.8323164747
$ws.Range("S6").Value = 0.1863055663917521
$ws.Range("T6").Value = 0.1863055663917521

# Row 7: updated NATMI metrics for this sending/target cluster pair
$ws.Range("A7").Value = "MuSCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.7137986666666666
$ws.Range("H7").Value = 2.141396
$ws.Range("I7").Value = 0.2157511749766052
$ws.Range("J7").Value = 0.2157511749766052
$ws.Range("O7").Value = 0.1364794819219223
$ws.Range("P7").Value = 0.1364794819219224
$ws.Range("Q7").Value = 0.1375442444088888
$ws.Range("R7").Value = 1.23789819968
$ws.Range("S7").Value = 0.02944560858485309
$ws.Range("T7").Value = 0.0294456085848531
